# Generate Report for handoff
# - The row describing file 78d9e4aa-...md moves to row 2 and keeps its
#   "Handed back: in sync with en-US" status.
# - The row describing file 26508feb-...md moves to row 3 and its status
#   changes to "Ready for handoff" with a refreshed handoff timestamp.
# Applies to the Overview sheet as well as the per-locale (zh-cn / de-de)
# detail sheets, including the hyperlinks that decorate the file-name /
# handoff-file / handback-file columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"

$ws1.Range("A3").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/e2e/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/e2e/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md"
$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-26 12:30:31"
$ws2.Range("E2").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md"
$ws2.Range("F2").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-01-26 12:29:19"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md"
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-01-26 12:30:31"
$ws2.Range("E3").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md"
$ws2.Range("F3").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf"
$ws2.Range("G3").Value = "2016-01-26 12:29:19"
$ws2.Range("H3").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/e2e/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5ca88405635ef0c3bb21fdd52673af46fcb9cee7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/05fef02f6be56dc738d88cad2136542ecec9edcc/e2e/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/601b964838776b0b272433c2c758816aa5a94328/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/e2e/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5ca88405635ef0c3bb21fdd52673af46fcb9cee7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/05fef02f6be56dc738d88cad2136542ecec9edcc/e2e/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/601b964838776b0b272433c2c758816aa5a94328/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md"
$ws3.Range("B2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf"
$ws3.Range("D2").Value = "2016-01-26 12:30:43"
$ws3.Range("E2").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md"
$ws3.Range("F2").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf"
$ws3.Range("G2").Value = "2016-01-26 12:29:39"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md"
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf"
$ws3.Range("D3").Value = "2016-01-26 12:30:43"
$ws3.Range("E3").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md"
$ws3.Range("F3").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf"
$ws3.Range("G3").Value = "2016-01-26 12:29:39"
$ws3.Range("H3").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/e2e/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d9f5d0492f4d5d912b973678a82e27543a2d191b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e67c4144042519038dd518a6b18fdab830a05744/e2e/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b29d735140eb965d8d21e791dcf322346526a9ac/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/e2e/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d9f5d0492f4d5d912b973678a82e27543a2d191b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e67c4144042519038dd518a6b18fdab830a05744/e2e/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b29d735140eb965d8d21e791dcf322346526a9ac/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/.localization-config", "", "", ".localization-config")
